$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 5).Value = 1
$ws.Cells.Item(2, 6).Value = 0.3333333333333333
$ws.Cells.Item(2, 7).Value = 0.164296
$ws.Cells.Item(2, 8).Value = 0.492888
$ws.Cells.Item(2, 9).Value = 0.04423686189757333
$ws.Cells.Item(2, 10).Value = 0.04423686189757334
$ws.Cells.Item(2, 11).Value = 2
$ws.Cells.Item(2, 12).Value = 0.6666666666666666
$ws.Cells.Item(2, 13).Value = 0.9029063333333333
$ws.Cells.Item(2, 14).Value = 2.708719
$ws.Cells.Item(2, 15).Value = 0.0007562739491265452
$ws.Cells.Item(2, 16).Value = 0.0007562739491265452
$ws.Cells.Item(2, 17).Value = 0.1483438989413333
$ws.Cells.Item(2, 18).Value = 1.335095090472
$ws.Cells.Item(2, 19).Value = 0.00003345518624424338
$ws.Cells.Item(2, 20).Value = 0.00003345518624424338

# Row 3
$ws.Cells.Item(3, 5).Value = 1
$ws.Cells.Item(3, 6).Value = 0.3333333333333333
$ws.Cells.Item(3, 7).Value = 0.164296
$ws.Cells.Item(3, 8).Value = 0.492888
$ws.Cells.Item(3, 9).Value = 0.04423686189757333
$ws.Cells.Item(3, 10).Value = 0.04423686189757334
$ws.Cells.Item(3, 13).Value = 63.06324899999999
$ws.Cells.Item(3, 14).Value = 189.189747
$ws.Cells.Item(3, 15).Value = 0.05282174972669441
$ws.Cells.Item(3, 16).Value = 0.0528217497266944
$ws.Cells.Item(3, 17).Value = 10.361039557704
$ws.Cells.Item(3, 18).Value = 93.24935601933599
$ws.Cells.Item(3, 19).Value = 0.002336668447847962
$ws.Cells.Item(3, 20).Value = 0.002336668447847962

# Row 4
$ws.Cells.Item(4, 5).Value = 1
$ws.Cells.Item(4, 6).Value = 0.3333333333333333
$ws.Cells.Item(4, 7).Value = 0.164296
$ws.Cells.Item(4, 8).Value = 0.492888
$ws.Cells.Item(4, 9).Value = 0.04423686189757333
$ws.Cells.Item(4, 10).Value = 0.04423686189757334
$ws.Cells.Item(4, 15).Value = 0.000007424203315745918
$ws.Cells.Item(4, 16).Value = 0.000007424203315745916
$ws.Cells.Item(4, 17).Value = 0.001456264978666667
$ws.Cells.Item(4, 18).Value = 0.013106384808
$ws.Cells.Item(4, 19).Value = 0.0000003284234567781582
$ws.Cells.Item(4, 20).Value = 0.0000003284234567781582

# Row 5
$ws.Cells.Item(5, 5).Value = 1
$ws.Cells.Item(5, 6).Value = 0.3333333333333333
$ws.Cells.Item(5, 7).Value = 0.164296
$ws.Cells.Item(5, 8).Value = 0.492888
$ws.Cells.Item(5, 9).Value = 0.04423686189757333
$ws.Cells.Item(5, 10).Value = 0.04423686189757334
$ws.Cells.Item(5, 13).Value = 1100.424519666667
$ws.Cells.Item(5, 14).Value = 3301.273559
$ws.Cells.Item(5, 15).Value = 0.9217150954425227
$ws.Cells.Item(5, 16).Value = 0.9217150954425226
$ws.Cells.Item(5, 17).Value = 180.7953468831547
$ws.Cells.Item(5, 18).Value = 1627.158121948392
$ws.Cells.Item(5, 19).Value = 0.0407737833859995
$ws.Cells.Item(5, 20).Value = 0.0407737833859995

# Row 6
$ws.Cells.Item(6, 5).Value = 1
$ws.Cells.Item(6, 6).Value = 0.3333333333333333
$ws.Cells.Item(6, 7).Value = 0.164296
$ws.Cells.Item(6, 8).Value = 0.492888
$ws.Cells.Item(6, 9).Value = 0.04423686189757333
$ws.Cells.Item(6, 10).Value = 0.04423686189757334
$ws.Cells.Item(6, 13).Value = 29.488383
$ws.Cells.Item(6, 14).Value = 88.465149
$ws.Cells.Item(6, 15).Value = 0.02469945667834066
$ws.Cells.Item(6, 16).Value = 0.02469945667834066
$ws.Cells.Item(6, 17).Value = 4.844823373368
$ws.Cells.Item(6, 18).Value = 43.603410360312
$ws.Cells.Item(6, 19).Value = 0.001092626454024851
$ws.Cells.Item(6, 20).Value = 0.001092626454024851

# Row 7
$ws.Cells.Item(7, 5).Value = 3
$ws.Cells.Item(7, 6).Value = 1
$ws.Cells.Item(7, 7).Value = 0.9811059999999999
$ws.Cells.Item(7, 8).Value = 2.943318
$ws.Cells.Item(7, 9).Value = 0.264163769226765
$ws.Cells.Item(7, 10).Value = 0.264163769226765
$ws.Cells.Item(7, 11).Value = 2
$ws.Cells.Item(7, 12).Value = 0.6666666666666666
$ws.Cells.Item(7, 13).Value = 0.9029063333333333
$ws.Cells.Item(7, 14).Value = 2.708719
$ws.Cells.Item(7, 15).Value = 0.0007562739491265452
$ws.Cells.Item(7, 16).Value = 0.0007562739491265452
$ws.Cells.Item(7, 17).Value = 0.8858468210713332
$ws.Cells.Item(7, 18).Value = 7.972621389641999
$ws.Cells.Item(7, 19).Value = 0.0001997801769692789
$ws.Cells.Item(7, 20).Value = 0.0001997801769692789

# Row 8
$ws.Cells.Item(8, 5).Value = 3
$ws.Cells.Item(8, 6).Value = 1
$ws.Cells.Item(8, 7).Value = 0.9811059999999999
$ws.Cells.Item(8, 8).Value = 2.943318
$ws.Cells.Item(8, 9).Value = 0.264163769226765
$ws.Cells.Item(8, 10).Value = 0.264163769226765
$ws.Cells.Item(8, 13).Value = 63.06324899999999
$ws.Cells.Item(8, 14).Value = 189.189747
$ws.Cells.Item(8, 15).Value = 0.05282174972669441
$ws.Cells.Item(8, 16).Value = 0.0528217497266944
$ws.Cells.Item(8, 17).Value = 61.87173197339398
$ws.Cells.Item(8, 18).Value = 556.8455877605459
$ws.Cells.Item(8, 19).Value = 0.01395359250495644
$ws.Cells.Item(8, 20).Value = 0.01395359250495644

# Row 9
$ws.Cells.Item(9, 5).Value = 3
$ws.Cells.Item(9, 6).Value = 1
$ws.Cells.Item(9, 7).Value = 0.9811059999999999
$ws.Cells.Item(9, 8).Value = 2.943318
$ws.Cells.Item(9, 9).Value = 0.264163769226765
$ws.Cells.Item(9, 10).Value = 0.264163769226765
$ws.Cells.Item(9, 15).Value = 0.000007424203315745918
$ws.Cells.Item(9, 16).Value = 0.000007424203315745916
$ws.Cells.Item(9, 17).Value = 0.008696196548666666
$ws.Cells.Item(9, 18).Value = 0.078265768938
$ws.Cells.Item(9, 19).Value = 0.000001961205531393288
$ws.Cells.Item(9, 20).Value = 0.000001961205531393288

# Row 10
$ws.Cells.Item(10, 5).Value = 3
$ws.Cells.Item(10, 6).Value = 1
$ws.Cells.Item(10, 7).Value = 0.9811059999999999
$ws.Cells.Item(10, 8).Value = 2.943318
$ws.Cells.Item(10, 9).Value = 0.264163769226765
$ws.Cells.Item(10, 10).Value = 0.264163769226765
$ws.Cells.Item(10, 13).Value = 1100.424519666667
$ws.Cells.Item(10, 14).Value = 3301.273559
$ws.Cells.Item(10, 15).Value = 0.9217150954425227
$ws.Cells.Item(10, 16).Value = 0.9217150954425226
$ws.Cells.Item(10, 17).Value = 1079.633098792084
$ws.Cells.Item(10, 18).Value = 9716.69788912876
$ws.Cells.Item(10, 19).Value = 0.2434837337653042
$ws.Cells.Item(10, 20).Value = 0.2434837337653042

# Row 11
$ws.Cells.Item(11, 5).Value = 3
$ws.Cells.Item(11, 6).Value = 1
$ws.Cells.Item(11, 7).Value = 0.9811059999999999
$ws.Cells.Item(11, 8).Value = 2.943318
$ws.Cells.Item(11, 9).Value = 0.264163769226765
$ws.Cells.Item(11, 10).Value = 0.264163769226765
$ws.Cells.Item(11, 13).Value = 29.488383
$ws.Cells.Item(11, 14).Value = 88.465149
$ws.Cells.Item(11, 15).Value = 0.02469945667834066
$ws.Cells.Item(11, 16).Value = 0.02469945667834066
$ws.Cells.Item(11, 17).Value = 28.931229491598
$ws.Cells.Item(11, 18).Value = 260.381065424382
$ws.Cells.Item(11, 19).Value = 0.006524701574003662
$ws.Cells.Item(11, 20).Value = 0.006524701574003661

# Row 12
$ws.Cells.Item(12, 7).Value = 2.382039
$ws.Cells.Item(12, 8).Value = 7.146117
$ws.Cells.Item(12, 9).Value = 0.6413663770124269
$ws.Cells.Item(12, 10).Value = 0.6413663770124269
$ws.Cells.Item(12, 11).Value = 2
$ws.Cells.Item(12, 12).Value = 0.6666666666666666
$ws.Cells.Item(12, 13).Value = 0.9029063333333333
$ws.Cells.Item(12, 14).Value = 2.708719
$ws.Cells.Item(12, 15).Value = 0.0007562739491265452
$ws.Cells.Item(12, 16).Value = 0.0007562739491265452
$ws.Cells.Item(12, 17).Value = 2.150758099347
$ws.Cells.Item(12, 18).Value = 19.356822894123
$ws.Cells.Item(12, 19).Value = 0.0004850486827801728
$ws.Cells.Item(12, 20).Value = 0.0004850486827801728

# Row 13
$ws.Cells.Item(13, 7).Value = 2.382039
$ws.Cells.Item(13, 8).Value = 7.146117
$ws.Cells.Item(13, 9).Value = 0.6413663770124269
$ws.Cells.Item(13, 10).Value = 0.6413663770124269
$ws.Cells.Item(13, 13).Value = 63.06324899999999
$ws.Cells.Item(13, 14).Value = 189.189747
$ws.Cells.Item(13, 15).Value = 0.05282174972669441
$ws.Cells.Item(13, 16).Value = 0.0528217497266944
$ws.Cells.Item(13, 17).Value = 150.219118584711
$ws.Cells.Item(13, 18).Value = 1351.972067262399
$ws.Cells.Item(13, 19).Value = 0.03387809424966715
$ws.Cells.Item(13, 20).Value = 0.03387809424966714

# Row 14
$ws.Cells.Item(14, 7).Value = 2.382039
$ws.Cells.Item(14, 8).Value = 7.146117
$ws.Cells.Item(14, 9).Value = 0.6413663770124269
$ws.Cells.Item(14, 10).Value = 0.6413663770124269
$ws.Cells.Item(14, 15).Value = 0.000007424203315745918
$ws.Cells.Item(14, 16).Value = 0.000007424203315745916
$ws.Cells.Item(14, 17).Value = 0.021113599683
$ws.Cells.Item(14, 18).Value = 0.190022397147
$ws.Cells.Item(14, 19).Value = 0.000004761634382823606
$ws.Cells.Item(14, 20).Value = 0.000004761634382823606

# Row 15
$ws.Cells.Item(15, 7).Value = 2.382039
$ws.Cells.Item(15, 8).Value = 7.146117
$ws.Cells.Item(15, 9).Value = 0.6413663770124269
$ws.Cells.Item(15, 10).Value = 0.6413663770124269
$ws.Cells.Item(15, 13).Value = 1100.424519666667
$ws.Cells.Item(15, 14).Value = 3301.273559
$ws.Cells.Item(15, 15).Value = 0.9217150954425227
$ws.Cells.Item(15, 16).Value = 0.9217150954425226
$ws.Cells.Item(15, 17).Value = 2621.254122402267
$ws.Cells.Item(15, 18).Value = 23591.2871016204
$ws.Cells.Item(15, 19).Value = 0.5911570714016341
$ws.Cells.Item(15, 20).Value = 0.591157071401634

# Row 16
$ws.Cells.Item(16, 7).Value = 2.382039
$ws.Cells.Item(16, 8).Value = 7.146117
$ws.Cells.Item(16, 9).Value = 0.6413663770124269
$ws.Cells.Item(16, 10).Value = 0.6413663770124269
$ws.Cells.Item(16, 13).Value = 29.488383
$ws.Cells.Item(16, 14).Value = 88.465149
$ws.Cells.Item(16, 15).Value = 0.02469945667834066
$ws.Cells.Item(16, 16).Value = 0.02469945667834066
$ws.Cells.Item(16, 17).Value = 70.242478352937
$ws.Cells.Item(16, 18).Value = 632.182305176433
$ws.Cells.Item(16, 19).Value = 0.01584140104396274
$ws.Cells.Item(16, 20).Value = 0.01584140104396274

# Row 17
$ws.Cells.Item(17, 7).Value = 0.1865656666666667
$ws.Cells.Item(17, 8).Value = 0.559697
$ws.Cells.Item(17, 9).Value = 0.05023299186323485
$ws.Cells.Item(17, 10).Value = 0.05023299186323486
$ws.Cells.Item(17, 11).Value = 2
$ws.Cells.Item(17, 12).Value = 0.6666666666666666
$ws.Cells.Item(17, 13).Value = 0.9029063333333333
$ws.Cells.Item(17, 14).Value = 2.708719
$ws.Cells.Item(17, 15).Value = 0.0007562739491265452
$ws.Cells.Item(17, 16).Value = 0.0007562739491265452
$ws.Cells.Item(17, 17).Value = 0.1684513220158889
$ws.Cells.Item(17, 18).Value = 1.516061898143
$ws.Cells.Item(17, 19).Value = 0.00003798990313285023
$ws.Cells.Item(17, 20).Value = 0.00003798990313285023

# Row 18
$ws.Cells.Item(18, 7).Value = 0.1865656666666667
$ws.Cells.Item(18, 8).Value = 0.559697
$ws.Cells.Item(18, 9).Value = 0.05023299186323485
$ws.Cells.Item(18, 10).Value = 0.05023299186323486
$ws.Cells.Item(18, 13).Value = 63.06324899999999
$ws.Cells.Item(18, 14).Value = 189.189747
$ws.Cells.Item(18, 15).Value = 0.05282174972669441
$ws.Cells.Item(18, 16).Value = 0.0528217497266944
$ws.Cells.Item(18, 17).Value = 11.765437091851
$ws.Cells.Item(18, 18).Value = 105.888933826659
$ws.Cells.Item(18, 19).Value = 0.002653394524222868
$ws.Cells.Item(18, 20).Value = 0.002653394524222868

# Row 19
$ws.Cells.Item(19, 7).Value = 0.1865656666666667
$ws.Cells.Item(19, 8).Value = 0.559697
$ws.Cells.Item(19, 9).Value = 0.05023299186323485
$ws.Cells.Item(19, 10).Value = 0.05023299186323486
$ws.Cells.Item(19, 15).Value = 0.000007424203315745918
$ws.Cells.Item(19, 16).Value = 0.000007424203315745916
$ws.Cells.Item(19, 17).Value = 0.001653655880777778
$ws.Cells.Item(19, 18).Value = 0.014882902927
$ws.Cells.Item(19, 19).Value = 0.0000003729399447508659
$ws.Cells.Item(19, 20).Value = 0.0000003729399447508659

# Row 20
$ws.Cells.Item(20, 7).Value = 0.1865656666666667
$ws.Cells.Item(20, 8).Value = 0.559697
$ws.Cells.Item(20, 9).Value = 0.05023299186323485
$ws.Cells.Item(20, 10).Value = 0.05023299186323486
$ws.Cells.Item(20, 13).Value = 1100.424519666667
$ws.Cells.Item(20, 14).Value = 3301.273559
$ws.Cells.Item(20, 15).Value = 0.9217150954425227
$ws.Cells.Item(20, 16).Value = 0.9217150954425226
$ws.Cells.Item(20, 17).Value = 205.3014341279581
$ws.Cells.Item(20, 18).Value = 1847.712907151623
$ws.Cells.Item(20, 19).Value = 0.04630050688958497
$ws.Cells.Item(20, 20).Value = 0.04630050688958498

# Row 21
$ws.Cells.Item(21, 7).Value = 0.1865656666666667
$ws.Cells.Item(21, 8).Value = 0.559697
$ws.Cells.Item(21, 9).Value = 0.05023299186323485
$ws.Cells.Item(21, 10).Value = 0.05023299186323486
$ws.Cells.Item(21, 13).Value = 29.488383
$ws.Cells.Item(21, 14).Value = 88.465149
$ws.Cells.Item(21, 15).Value = 0.02469945667834066
$ws.Cells.Item(21, 16).Value = 0.02469945667834066
$ws.Cells.Item(21, 17).Value = 5.501519833316999
$ws.Cells.Item(21, 18).Value = 28.931229491598
$ws.Cells.Item(21, 19).Value = 0.001240727606349408
$ws.Cells.Item(21, 20).Value = 0.001240727606349408
